# Fruta / hortaliza, semanal
# Insert two new weekly rows (date 2022-12-23 / serial 44918) right above the
# existing row 552 block, pushing the rest of the price-history rows down by
# two (old row N -> new row N+2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 552 (each Insert() pushes row 552.. down by one).
$ws.Rows.Item(552).Insert()
$ws.Rows.Item(552).Insert()

# New row 552: Copenhague -> Crespo record, Primera, 500 / 2200 / 2200 / 2200
$ws.Cells.Item(552, 1).Value = 4
$ws.Cells.Item(552, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(552, 3).Value = "Los Lagos"
$ws.Cells.Item(552, 4).Value = 44918
$ws.Cells.Item(552, 5).Value = 10
$ws.Cells.Item(552, 6).Value = 100112006
$ws.Cells.Item(552, 7).Value = "Repollo"
$ws.Cells.Item(552, 8).Value = "Crespo record"
$ws.Cells.Item(552, 9).Value = "Primera"
$ws.Cells.Item(552, 10).Value = 500
$ws.Cells.Item(552, 11).Value = 2200
$ws.Cells.Item(552, 12).Value = 2200
$ws.Cells.Item(552, 13).Value = 2200
$ws.Cells.Item(552, 14).Value = "$/unidad"
$ws.Cells.Item(552, 15).Value = "Región Metropolitana"
$ws.Cells.Item(552, 16).Value = 2200
$ws.Cells.Item(552, 17).Value = 1
$ws.Cells.Item(552, 18).Value = "Hortaliza"

# New row 553: Crespo record, Primera -> Segunda, 500 / 2000 / 2000 / 2000
$ws.Cells.Item(553, 1).Value = 4
$ws.Cells.Item(553, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(553, 3).Value = "Los Lagos"
$ws.Cells.Item(553, 4).Value = 44918
$ws.Cells.Item(553, 5).Value = 10
$ws.Cells.Item(553, 6).Value = 100112006
$ws.Cells.Item(553, 7).Value = "Repollo"
$ws.Cells.Item(553, 8).Value = "Crespo record"
$ws.Cells.Item(553, 9).Value = "Segunda"
$ws.Cells.Item(553, 10).Value = 500
$ws.Cells.Item(553, 11).Value = 2000
$ws.Cells.Item(553, 12).Value = 2000
$ws.Cells.Item(553, 13).Value = 2000
$ws.Cells.Item(553, 14).Value = "$/unidad"
$ws.Cells.Item(553, 15).Value = "Región Metropolitana"
$ws.Cells.Item(553, 16).Value = 2000
$ws.Cells.Item(553, 17).Value = 1
$ws.Cells.Item(553, 18).Value = "Hortaliza"
